# Commit: Add PF/1.0.5 to meta-sheet
# Append a new row (row 3) to the meta-sheet with a new package version
# (PF/1.0.5) and mark it unsupported ("X") for the dev2/sit2/uat2 columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "PF/1.0.5"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"

# New row should use the default "Normal" cell style (no explicit
# formatting carried over from the column), matching the rest of the
# sheet's existing rows.
$ws.Range("A3:D3").Style = "Normal"
